$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")
$wsTransactions = $wb.Worksheets.Item("Transactions")

# Insert a new blank column before column N on the "Repayment Schedule" sheet.
$wsSchedule.Columns("N:N").Insert()

# Update sheet selections/active states.
$wsTransactions.Range("B8").Select()
$wsSchedule.Activate()
$wsSchedule.Range("S5").Select()
